# CTF-101 H4ck3r k1d3rgart3n.pptx — "Final edition of slides presented. Added a pdf version"
#
# The diff reduces to:
#   1. Delete slide 14 ("Warm up #1") — all later slides shift up by one.
#   2. On what is now the final slide (previously slide 28, "That's all there
#      is to it…"):
#        - shorten the title to "That's it" (2nd paragraph "Go have fun now"
#          is unchanged)
#        - nudge the "Miki Demeter…" textbox up slightly (top: 53.35pt -> 26.35pt)
#        - nudge the "Stacy Watts…" textbox down slightly (top: 533.35pt -> 558.35pt)

$p = $ppt.ActivePresentation

# 1. Remove "Warm up #1" slide.
$warmupIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.Count -gt 0 -and $candidate.Shapes.Item(1).Name -eq "Warm up #1") {
        $warmupIndex = $i
        break
    }
}
if ($warmupIndex -eq -1) {
    $warmupIndex = 14
}
$p.Slides.Item($warmupIndex).Delete()

# 2. Tweak the closing slide (now the last slide in the deck).
$closing = $p.Slides.Item($p.Slides.Count)
$title = $closing.Shapes.Item(1)
$miki = $closing.Shapes.Item(2)
$stacy = $closing.Shapes.Item(3)

$rsquo = [char]0x2019
$hellip = [char]0x2026

$title.TextFrame.TextRange.Text = "That" + $rsquo + "s it" + [char]0x0B + "Go have fun now"
$title.Name = "That" + $rsquo + "s it" + $hellip

$miki.Top = 334620 / 12700.0
$stacy.Top = 7091020 / 12700.0
